$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 35
$ws.Range("I2").Value = 66
$ws.Range("J2").Value = 326
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 98
$ws.Range("M2").Value = 6
$ws.Range("N2").Value = 58
$ws.Range("P2").Value = 2
$ws.Range("R2").Value = 4
$ws.Range("S2").Value = 28
$ws.Range("T2").Value = 57
$ws.Range("V2").Value = 505
$ws.Range("X2").Value = 509
$ws.Range("Z2").Value = 8
$ws.Range("AA2").Value = 5
